# Generate Report for Handoff
# Row 3 ("b.md") has been handed off again: update its status and the
# handoff file/datetime details on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Status column for b.md (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row 3 (b.md) handoff details ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-03-10 09:20:14"

# --- de-de sheet: row 3 (b.md) handoff details ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-03-10 09:20:25"
